# The source data sheet got its weekly refresh: a new week (2 price rows,
# one per quality grade) was inserted right above what used to be row 474,
# pushing the rest of the table (old rows 474:580) down to 476:582.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(474).Resize(2).EntireRow.Insert()

# Populate the two newly inserted rows (474 and 475) with the new week's data.
$ws.Range("A474").Value = 6
$ws.Range("B474").Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range("C474").Value = 'Metropolitana'
$ws.Range("D474").Value = 44511
$ws.Range("E474").Value = 13
$ws.Range("F474").Value = 100112040
$ws.Range("G474").Value = 'Cilantro'
$ws.Range("H474").Value = 'Sin especificar'
$ws.Range("I474").Value = 'Primera'
$ws.Range("J474").Value = 630
$ws.Range("K474").Value = 4500
$ws.Range("L474").Value = 5000
$ws.Range("M474").Value = 4722
$ws.Range("N474").Value = '$/caja 36 atados'
$ws.Range("O474").Value = 'Región Metropolitana'
$ws.Range("P474").Value = 131
$ws.Range("Q474").Value = 36
$ws.Range("R474").Value = 'Hortaliza'

$ws.Range("A475").Value = 6
$ws.Range("B475").Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range("C475").Value = 'Metropolitana'
$ws.Range("D475").Value = 44511
$ws.Range("E475").Value = 13
$ws.Range("F475").Value = 100112040
$ws.Range("G475").Value = 'Cilantro'
$ws.Range("H475").Value = 'Sin especificar'
$ws.Range("I475").Value = 'Primera'
$ws.Range("J475").Value = 270
$ws.Range("K475").Value = 11000
$ws.Range("L475").Value = 12000
$ws.Range("M475").Value = 11407
$ws.Range("N475").Value = '$/docena de atados'
$ws.Range("O475").Value = 'Región Metropolitana'
$ws.Range("P475").Value = 3802
$ws.Range("Q475").Value = 3
$ws.Range("R475").Value = 'Hortaliza'
